$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2)
}

# Student copy (first occurrence block)
Replace-Text "John Vincent" "asdf"
Replace-Text "12-ambot" "adsf"
Replace-Text "example teacher" "asdf"
Replace-Text "2024-03-27   " "2024-03-20   "
Replace-Text "  09:01:00   " "  19:08:00   "
Replace-Text "19:57:00" "08:09:00"
Replace-Text "                                 Example Counselor________________ " "                                 asdfasdf________________ "
Replace-Text "                          example teacher_________ " "                          asdf_________ "

# Teacher copy (second occurrence block)
Replace-Text "John Vincent " "asdf "
Replace-Text "12-ambot" "adsf"
Replace-Text "example teacher" "asdf"
Replace-Text "          2024-03-27   " "          2024-03-20   "
Replace-Text "  09:01:00   " "  19:08:00   "
Replace-Text "19:57:00" "08:09:00"
Replace-Text "                                 Example Counselor________________ " "                                 asdfasdf________________ "
Replace-Text "                          example teacher_________ " "                          asdf_________ "
